$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.574.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.56%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.102.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +9.84%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.658"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.42%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.85"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.14%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.95"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.40%  "

$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.379"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0744"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.50%  "

$ws.Range("E12").Value = "  +0.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.408.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.837"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.099.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +9.63%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.29%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.568.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.86%  "

$ws.Range("E20").Value = "  -3.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "240.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.62%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.10%  "

$ws.Range("E24").Value = "  +0.15%  "

$ws.Range("E25").Value = "  -4.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -10.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +56.05%  "

$ws.Range("E31").Value = "  -5.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.99%  "

$ws.Range("E33").Value = "  -0.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +23.78%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.984"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +12.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0899"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.26%  "

$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("E38").Value = "  -1.93%  "

$ws.Range("E39").Value = "  -5.74%  "

$ws.Range("E40").Value = "  -11.55%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0224"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.67%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "98.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.75%  "

$ws.Range("E44").Value = "  -3.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.339.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0845"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.75%  "

$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.79%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.297.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.81%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.35%  "
